$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholders: "20/01/2022" -> "21/01/2022" on the slide master and
#    on every custom (slide) layout that carries a "Date Placeholder *" shape.
# ---------------------------------------------------------------------------
$newDate = "21/01/2022"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------------
# 2) Slide 1 textbox ("TextBox 129"): reposition/resize and update its text,
#    keeping the existing two-run split (the second run keeps its err="1"
#    spell-check flag) by editing each run's characters in place.
# ---------------------------------------------------------------------------
$slide = $p.Slides.Item(1)
$textBox = $null
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 129") {
        $textBox = $sh
        break
    }
}

$textBox.Left = 342.04582677165354
$textBox.Top = 473.2844094488189
$textBox.Width = 137.95417322834646
$textBox.Height = 21.810944881889764

$tr = $textBox.TextFrame2.TextRange
$firstLen = "Memory issues when getting ".Length
$tr.Characters(1, $firstLen).Text = "These file are "

$newFirstLen = "These file are ".Length
$totalLen = $tr.Length
$tr.Characters($newFirstLen + 1, $totalLen - $newFirstLen).Text = "multiframe"
